# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets per the scheduled
# market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 236
$ws.Range("I55").Value = 98.333336
$ws.Range("K55").Value = 98.333336
$ws.Range("M55").Value = 115.666664

# Row 70
$ws.Range("H70").Value = 52318.734
$ws.Range("I70").Value = 4464.2
$ws.Range("J70").Value = 76246
$ws.Range("K70").Value = 13392.6
$ws.Range("L70").Value = 228738
$ws.Range("M70").Value = -13122.6
$ws.Range("N70").Value = -229278

# Row 73
$ws.Range("H73").Value = 52318.734
$ws.Range("I73").Value = 4464.2
$ws.Range("J73").Value = 76246
$ws.Range("K73").Value = 13392.6
$ws.Range("L73").Value = 228738
$ws.Range("M73").Value = -12456.6
$ws.Range("N73").Value = -230610

# Row 74
$ws.Range("H74").Value = 4499.5
$ws.Range("I74").Value = 4499.5
$ws.Range("K74").Value = 4499.5
$ws.Range("M74").Value = -3563.5

# Row 77
$ws.Range("H77").Value = 4499.5
$ws.Range("I77").Value = 4499.5
$ws.Range("K77").Value = 22497.5
$ws.Range("M77").Value = -17817.5

# Row 106
$ws.Range("H106").Value = 31595.273
$ws.Range("I106").Value = 33838.668
$ws.Range("K106").Value = 33838.668
$ws.Range("M106").Value = -33207.668

# Row 111
$ws.Range("H111").Value = 3716.5715
$ws.Range("I111").Value = 3836.8333
$ws.Range("J111").Value = 2995
$ws.Range("K111").Value = 11510.4999
$ws.Range("L111").Value = 8985
$ws.Range("M111").Value = -8443.499899999999
$ws.Range("N111").Value = -15119

# Row 113
$ws.Range("H113").Value = 3598.6
$ws.Range("I113").Value = 3747
$ws.Range("J113").Value = 3499.6667
$ws.Range("K113").Value = 3747
$ws.Range("L113").Value = 3499.6667
$ws.Range("M113").Value = -493
$ws.Range("N113").Value = -10007.6667

# Row 116
$ws.Range("H116").Value = 5350
$ws.Range("I116").Value = 4937.5
$ws.Range("K116").Value = 4937.5
$ws.Range("M116").Value = -1495.5

# Row 137
$ws.Range("H137").Value = 3163.88
$ws.Range("I137").Value = 1461.3636
$ws.Range("K137").Value = 4384.0908
$ws.Range("M137").Value = -1834.0908

# Row 141
$ws.Range("H141").Value = 5439.5454
$ws.Range("I141").Value = 4783.5
$ws.Range("K141").Value = 14350.5
$ws.Range("M141").Value = -9170.5

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 5465.533
$ws.Range("I63").Value = 3183.1667
$ws.Range("J63").Value = 6987.1113
$ws.Range("K63").Value = 3183.1667
$ws.Range("L63").Value = 6987.1113
$ws.Range("M63").Value = -2497.1667
$ws.Range("N63").Value = -8359.1113

# Row 66
$ws.Range("H66").Value = 5465.533
$ws.Range("I66").Value = 3183.1667
$ws.Range("J66").Value = 6987.1113
$ws.Range("K66").Value = 15915.8335
$ws.Range("L66").Value = 34935.5565
$ws.Range("M66").Value = -12483.8335
$ws.Range("N66").Value = -41799.5565

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 131
$ws.Range("H131").Value = 200357.5
$ws.Range("J131").Value = 200357.5
$ws.Range("L131").Value = 200357.5
$ws.Range("N131").Value = -210437.5

# Row 132
$ws.Range("H132").Value = 1917.9
$ws.Range("I132").Value = 1612.5186
$ws.Range("K132").Value = 4837.5558
$ws.Range("M132").Value = -2307.5558

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4016.75
$ws.Range("I134").Value = 3624.4
$ws.Range("K134").Value = 10873.2
$ws.Range("M134").Value = -8338.200000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4250.1924
$ws.Range("I31").Value = 2646.0557
$ws.Range("J31").Value = 7859.5
$ws.Range("K31").Value = 2646.0557
$ws.Range("L31").Value = 7859.5
$ws.Range("M31").Value = -2351.0557
$ws.Range("N31").Value = -8449.5

# Row 34
$ws.Range("H34").Value = 4250.1924
$ws.Range("I34").Value = 2646.0557
$ws.Range("J34").Value = 7859.5
$ws.Range("K34").Value = 2646.0557
$ws.Range("L34").Value = 7859.5
$ws.Range("M34").Value = -2444.0557
$ws.Range("N34").Value = -8263.5

# Row 62
$ws.Range("H62").Value = 73898.664
$ws.Range("J62").Value = 73898.664
$ws.Range("L62").Value = 73898.664
$ws.Range("N62").Value = -75146.664

# Row 65
$ws.Range("H65").Value = 73898.664
$ws.Range("J65").Value = 73898.664
$ws.Range("L65").Value = 369493.32
$ws.Range("N65").Value = -375733.32

# Row 107
$ws.Range("H107").Value = 665.7059
$ws.Range("I107").Value = 404.2
$ws.Range("J107").Value = 1039.2858
$ws.Range("K107").Value = 404.2
$ws.Range("L107").Value = 1039.2858
$ws.Range("M107").Value = 1515.8
$ws.Range("N107").Value = -4879.2858

# Row 119
$ws.Range("H119").Value = 39999
$ws.Range("J119").Value = 39999
$ws.Range("L119").Value = 39999
$ws.Range("N119").Value = -49675

# Row 132
$ws.Range("H132").Value = 3752.875
$ws.Range("I132").Value = 2337.6667
$ws.Range("K132").Value = 7013.000100000001
$ws.Range("M132").Value = -4483.000100000001

# Row 138
$ws.Range("H138").Value = 99166.664
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

$ws = $wb.Worksheets.Item("CUL")
# Row 35
$ws.Range("H35").Value = 107.5
$ws.Range("J35").Value = 107.5
$ws.Range("L35").Value = 322.5
$ws.Range("N35").Value = -898.5

# Row 36
$ws.Range("H36").Value = 1766.5555
$ws.Range("I36").Value = 2071.4285
$ws.Range("J36").Value = 699.5
$ws.Range("K36").Value = 6214.2855
$ws.Range("L36").Value = 2098.5
$ws.Range("M36").Value = -6045.2855
$ws.Range("N36").Value = -2436.5

# Row 81
$ws.Range("H81").Value = 2711.7144
$ws.Range("J81").Value = 2711.7144
$ws.Range("L81").Value = 8135.1432
$ws.Range("N81").Value = -10381.1432

# Row 84
$ws.Range("H84").Value = 2711.7144
$ws.Range("J84").Value = 2711.7144
$ws.Range("L84").Value = 24405.4296
$ws.Range("N84").Value = -35637.4296

# Row 103
$ws.Range("H103").Value = 4333.3335
$ws.Range("J103").Value = 4333.3335
$ws.Range("L103").Value = 13000.0005
$ws.Range("N103").Value = -14758.0005

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3897.6428
$ws.Range("I126").Value = 2542.3333
$ws.Range("J126").Value = 4914.125
$ws.Range("K126").Value = 7626.999899999999
$ws.Range("L126").Value = 14742.375
$ws.Range("M126").Value = -5156.999899999999
$ws.Range("N126").Value = -19682.375

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 942.25
$ws.Range("J46").Value = 908
$ws.Range("L46").Value = 908
$ws.Range("N46").Value = -1284

# Row 93
$ws.Range("H93").Value = 3097.2
$ws.Range("I93").Value = 3097.2
$ws.Range("K93").Value = 3097.2
$ws.Range("M93").Value = -1849.2

# Row 134
$ws.Range("H134").Value = 110000
$ws.Range("J134").Value = 110000
$ws.Range("L134").Value = 110000
$ws.Range("N134").Value = -120140

# Row 136
$ws.Range("H136").Value = 3250
$ws.Range("I136").Value = 3250
$ws.Range("K136").Value = 9750
$ws.Range("M136").Value = -7200

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 99998
$ws.Range("J46").Value = 99998
$ws.Range("L46").Value = 99998
$ws.Range("N46").Value = -100460

# Row 126
$ws.Range("H126").Value = 2092.3157
$ws.Range("J126").Value = 5755
$ws.Range("L126").Value = 17265
$ws.Range("N126").Value = -22205

# Row 132
$ws.Range("H132").Value = 2125.375
$ws.Range("I132").Value = 1857.5714
$ws.Range("K132").Value = 5572.7142
$ws.Range("M132").Value = -3042.7142

# Row 134
$ws.Range("H134").Value = 99998
$ws.Range("J134").Value = 99998
$ws.Range("L134").Value = 299994
$ws.Range("N134").Value = -305064

# Row 136
$ws.Range("H136").Value = 6507.222
$ws.Range("J136").Value = 6331.3335
$ws.Range("L136").Value = 18994.0005
$ws.Range("N136").Value = -24094.0005
